# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.585.94'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.40%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.649.74'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.77%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.14'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.20'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.63%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.648.56'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.69%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.61'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.153'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.47'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.122.81'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.472.78'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000145'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.624.88'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.43'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '340.71'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.36'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.73'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.35%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.20'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.68'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +6.32%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +5.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.166'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '548.26'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +18.61%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.44'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.79'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.82'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +14.21%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.82%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '175.36'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.98%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.92'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +9.49%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.09'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +7.44%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '170.17'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +7.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.25'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.74'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.45'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +6.69%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0556'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +3.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0960'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.38%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.72'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.71'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.18%  '
